$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Correct the "Integer min" value for rule R20 (row 10) from 18 to 1.
$ws.Range("C10").Value = 1
